# Duplication formats 21 > 22 (pour janvier-fevrier avant maj fg/genrsa)
# Adds rows 100:107 to Feuil1, duplicating rows 92:99 (the "2021" rsa/rum
# blocks) but with the year column (B) updated to 2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Copy the source block (rows 92:99, cols A:E) into the new block (rows
# 100:107) - this carries over cell values, string vs numeric types and
# styles exactly like the existing rows.
$src = $ws.Range("A92:E99")
$dst = $ws.Range("A100:E107")
$src.Copy($dst)

# The duplicated block still says "2021" in column B; bump it to 2022.
$ws.Range("B100:B107").Value = 2022

# Move the view / selection the same way the diff shows (scrolled down to
# the newly added rows, active cell on the last rsa row of the new block).
$win = $excel.ActiveWindow
$win.ScrollRow = 87
$ws.Range("B99").Select()
